$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 12615.6
$ws.Range("I40").Value = 14870.875
$ws.Range("J40").Value = 3594.5
$ws.Range("K40").Value = 14870.875
$ws.Range("L40").Value = 3594.5
$ws.Range("M40").Value = -14695.875
$ws.Range("N40").Value = -3944.5
$ws.Range("H74").Value = 7556.4165
$ws.Range("I74").Value = 6136.6
$ws.Range("J74").Value = 8570.571
$ws.Range("K74").Value = 6136.6
$ws.Range("L74").Value = 8570.571
$ws.Range("M74").Value = -5200.6
$ws.Range("N74").Value = -10442.571
$ws.Range("H77").Value = 7556.4165
$ws.Range("I77").Value = 6136.6
$ws.Range("J77").Value = 8570.571
$ws.Range("K77").Value = 30683
$ws.Range("L77").Value = 42852.855
$ws.Range("M77").Value = -26003
$ws.Range("N77").Value = -52212.855
$ws.Range("H86").Value = 23811462
$ws.Range("I86").Value = 24693364
$ws.Range("J86").Value = 18520052
$ws.Range("K86").Value = 24693364
$ws.Range("L86").Value = 18520052
$ws.Range("M86").Value = -24692241
$ws.Range("N86").Value = -18522298
$ws.Range("H89").Value = 23811462
$ws.Range("I89").Value = 24693364
$ws.Range("J89").Value = 18520052
$ws.Range("K89").Value = 123466820
$ws.Range("L89").Value = 92600260
$ws.Range("M89").Value = -123461204
$ws.Range("N89").Value = -92611492
$ws.Range("H100").Value = 2455.4167
$ws.Range("J100").Value = 3160.8333
$ws.Range("L100").Value = 3160.8333
$ws.Range("N100").Value = -4242.8333
$ws.Range("H112").Value = 6028.222
$ws.Range("J112").Value = 6171.914
$ws.Range("L112").Value = 18515.742
$ws.Range("N112").Value = -20731.742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 105000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H132").Value = 5958.163
$ws.Range("I132").Value = 4445.657
$ws.Range("K132").Value = 13336.971
$ws.Range("M132").Value = -10806.971

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6175037.5
$ws.Range("I20").Value = 8334849.5
$ws.Range("K20").Value = 8334849.5
$ws.Range("M20").Value = -8334602.5
$ws.Range("H40").Value = 53942
$ws.Range("J40").Value = 53942
$ws.Range("L40").Value = 53942
$ws.Range("N40").Value = -54472
$ws.Range("H94").Value = 1020.10345
$ws.Range("I94").Value = 389.1905
$ws.Range("K94").Value = 389.1905
$ws.Range("M94").Value = 61.80950000000001
$ws.Range("I107").Value = 53573460
$ws.Range("J107").Value = 5747
$ws.Range("K107").Value = 53573460
$ws.Range("L107").Value = 5747
$ws.Range("M107").Value = -53571540
$ws.Range("N107").Value = -9587
$ws.Range("H134").Value = 4739.197
$ws.Range("I134").Value = 2861.86
$ws.Range("K134").Value = 8585.58
$ws.Range("M134").Value = -6050.58

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2150
$ws.Range("I107").Value = 1320.2
$ws.Range("J107").Value = 3187.25
$ws.Range("K107").Value = 1320.2
$ws.Range("L107").Value = 3187.25
$ws.Range("M107").Value = 599.8
$ws.Range("N107").Value = -7027.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 70.36364
$ws.Range("I40").Value = 47.4
$ws.Range("J40").Value = 89.5
$ws.Range("K40").Value = 189.6
$ws.Range("L40").Value = 358
$ws.Range("M40").Value = -120.6
$ws.Range("N40").Value = -496
$ws.Range("H116").Value = 2366
$ws.Range("H132").Value = 19832.166
$ws.Range("J132").Value = 22748.5
$ws.Range("L132").Value = 204736.5
$ws.Range("N132").Value = -209796.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6373.3257
$ws.Range("I70").Value = 5111.8623
$ws.Range("J70").Value = 8986.357
$ws.Range("K70").Value = 5111.8623
$ws.Range("L70").Value = 8986.357
$ws.Range("M70").Value = -4841.8623
$ws.Range("N70").Value = -9526.357
$ws.Range("H73").Value = 6373.3257
$ws.Range("I73").Value = 5111.8623
$ws.Range("J73").Value = 8986.357
$ws.Range("K73").Value = 5111.8623
$ws.Range("L73").Value = 8986.357
$ws.Range("M73").Value = -4175.8623
$ws.Range("N73").Value = -10858.357
$ws.Range("H80").Value = 168439
$ws.Range("I80").Value = 1227.5
$ws.Range("J80").Value = 252044.75
$ws.Range("K80").Value = 1227.5
$ws.Range("L80").Value = 252044.75
$ws.Range("M80").Value = -229.5
$ws.Range("N80").Value = -254040.75
$ws.Range("H83").Value = 168439
$ws.Range("I83").Value = 1227.5
$ws.Range("J83").Value = 252044.75
$ws.Range("K83").Value = 6137.5
$ws.Range("L83").Value = 1260223.75
$ws.Range("M83").Value = -1145.5
$ws.Range("N83").Value = -1270207.75
$ws.Range("H97").Value = 2436.0356
$ws.Range("J97").Value = 2733.25
$ws.Range("L97").Value = 2733.25
$ws.Range("N97").Value = -3725.25
$ws.Range("I113").Value = 3935.2856
$ws.Range("J113").Value = 6951.875
$ws.Range("K113").Value = 3935.2856
$ws.Range("L113").Value = 6951.875
$ws.Range("M113").Value = -1765.2856
$ws.Range("N113").Value = -11291.875
$ws.Range("H122").Value = 55391.6
$ws.Range("I122").Value = 170065.67
$ws.Range("K122").Value = 510197.01
$ws.Range("M122").Value = -507747.01

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 12347278
$ws.Range("I46").Value = 849
$ws.Range("J46").Value = 15874830
$ws.Range("K46").Value = 849
$ws.Range("L46").Value = 15874830
$ws.Range("M46").Value = -661
$ws.Range("N46").Value = -15875206
$ws.Range("H55").Value = 52631936
$ws.Range("I55").Value = 200000060
$ws.Range("J55").Value = 463.2143
$ws.Range("K55").Value = 200000060
$ws.Range("L55").Value = 463.2143
$ws.Range("M55").Value = -199999887
$ws.Range("N55").Value = -809.2143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5119.933
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 5599.6665
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5599.6665
$ws.Range("M4").Value = -4887
$ws.Range("N4").Value = -5825.6665
$ws.Range("H132").Value = 19255874
$ws.Range("I132").Value = 21745710
$ws.Range("K132").Value = 65237130
$ws.Range("M132").Value = -65234600

Write-Host "Applied 159 cell updates across 8 sheets"
